$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 (F column quantity changes)
$ws.Range("F2").Value = -192
$ws.Range("F3").Value = -406
$ws.Range("F4").Value = -660

# Row 5: location DC_002 -> DC_001, quantity -33 -> -79, horizon_days 1 -> 4
$ws.Range("B5").Value = "DC_001"
$ws.Range("F5").Value = -79
$ws.Range("H5").Value = 4

# New row 6: MAT_B / DC_002 / ... / -33 / ... / 1
$ws.Range("A6").Value = "MAT_B"
$ws.Range("B6").Value = "DC_002"
$ws.Range("C6").Value = 45294
$ws.Range("D6").Value = "Distribution Demand - Forecast"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -33
$ws.Range("G6").Value = 45293
$ws.Range("H6").Value = 1

# New row 7: MAT_B / PLANT_001 / ... / -79 / ... / 1
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = 45294
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -79
$ws.Range("G7").Value = 45293
$ws.Range("H7").Value = 1

# Copy formatting (number format/style) from existing date cells so C6/C7/G6/G7 match style s="2"
$ws.Range("C5").Copy()
$ws.Range("C6:C7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G5").Copy()
$ws.Range("G6:G7").PasteSpecial(-4122)  # xlPasteFormats
